$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 445.08334
$ws.Range("I33").Value = 427.22223
$ws.Range("K33").Value = 427.22223
$ws.Range("M33").Value = -198.22223
$ws.Range("H76").Value = 12351655
$ws.Range("I76").Value = 18522852
$ws.Range("K76").Value = 18522852
$ws.Range("M76").Value = -18522537
$ws.Range("H79").Value = 12351655
$ws.Range("I79").Value = 18522852
$ws.Range("K79").Value = 18522852
$ws.Range("M79").Value = -18521760
$ws.Range("H93").Value = 99499.5
$ws.Range("J93").Value = 99499.5
$ws.Range("L93").Value = 99499.5
$ws.Range("N93").Value = -104491.5
$ws.Range("H132").Value = 4286.2085
$ws.Range("I132").Value = 2369.2
$ws.Range("J132").Value = 13871.25
$ws.Range("K132").Value = 7107.599999999999
$ws.Range("L132").Value = 41613.75
$ws.Range("M132").Value = -4577.599999999999
$ws.Range("N132").Value = -46673.75
$ws.Range("H135").Value = 2652.2979
$ws.Range("I135").Value = 2104.1482
$ws.Range("J135").Value = 3392.3
$ws.Range("K135").Value = 18937.3338
$ws.Range("L135").Value = 30530.7
$ws.Range("M135").Value = -16402.3338
$ws.Range("N135").Value = -35600.7
$ws.Range("H137").Value = 1717.4412
$ws.Range("I137").Value = 1441.7084
$ws.Range("J137").Value = 2379.2
$ws.Range("K137").Value = 4325.1252
$ws.Range("L137").Value = 7137.599999999999
$ws.Range("M137").Value = -1775.1252
$ws.Range("N137").Value = -12237.6
$ws.Range("H141").Value = 18525132
$ws.Range("I141").Value = 21743880
$ws.Range("K141").Value = 65231640
$ws.Range("M141").Value = -65226460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 21415.834
$ws.Range("I28").Value = 6713
$ws.Range("J28").Value = 41999.8
$ws.Range("K28").Value = 6713
$ws.Range("L28").Value = 41999.8
$ws.Range("M28").Value = -6521
$ws.Range("N28").Value = -42383.8
$ws.Range("H45").Value = 2717.9167
$ws.Range("I45").Value = 2313.6
$ws.Range("K45").Value = 2313.6
$ws.Range("M45").Value = -1936.6
$ws.Range("H62").Value = 83414.336
$ws.Range("J62").Value = 83414.336
$ws.Range("L62").Value = 83414.336
$ws.Range("N62").Value = -84662.336
$ws.Range("H65").Value = 83414.336
$ws.Range("J65").Value = 83414.336
$ws.Range("L65").Value = 250243.008
$ws.Range("N65").Value = -256483.008
$ws.Range("H88").Value = 3563
$ws.Range("I88").Value = 2750.1667
$ws.Range("K88").Value = 2750.1667
$ws.Range("M88").Value = -2344.1667
$ws.Range("H91").Value = 3563
$ws.Range("I91").Value = 2750.1667
$ws.Range("K91").Value = 2750.1667
$ws.Range("M91").Value = -1346.1667
$ws.Range("H99").Value = 21415.834
$ws.Range("I99").Value = 6713
$ws.Range("J99").Value = 41999.8
$ws.Range("K99").Value = 6713
$ws.Range("L99").Value = 41999.8
$ws.Range("M99").Value = -3718
$ws.Range("N99").Value = -47989.8
$ws.Range("H102").Value = 27780036
$ws.Range("I102").Value = 45456370
$ws.Range("K102").Value = 45456370
$ws.Range("M102").Value = -45454748
$ws.Range("H122").Value = 4006.8484
$ws.Range("I122").Value = 3194.5938
$ws.Range("K122").Value = 9583.7814
$ws.Range("M122").Value = -7133.7814
$ws.Range("H138").Value = 65000
$ws.Range("I138").Value = 30000
$ws.Range("K138").Value = 30000
$ws.Range("M138").Value = -24860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2750.7058
$ws.Range("I80").Value = 1725.875
$ws.Range("J80").Value = 3661.6667
$ws.Range("K80").Value = 1725.875
$ws.Range("L80").Value = 3661.6667
$ws.Range("M80").Value = -727.875
$ws.Range("N80").Value = -5657.6667
$ws.Range("H83").Value = 2750.7058
$ws.Range("I83").Value = 1725.875
$ws.Range("J83").Value = 3661.6667
$ws.Range("K83").Value = 8629.375
$ws.Range("L83").Value = 18308.3335
$ws.Range("M83").Value = -3637.375
$ws.Range("N83").Value = -28292.3335
$ws.Range("H86").Value = 2718.8965
$ws.Range("I86").Value = 1205.5883
$ws.Range("J86").Value = 4862.75
$ws.Range("K86").Value = 1205.5883
$ws.Range("L86").Value = 4862.75
$ws.Range("M86").Value = -82.58829999999989
$ws.Range("N86").Value = -7108.75
$ws.Range("H89").Value = 2718.8965
$ws.Range("I89").Value = 1205.5883
$ws.Range("J89").Value = 4862.75
$ws.Range("K89").Value = 6027.941499999999
$ws.Range("L89").Value = 24313.75
$ws.Range("M89").Value = -411.941499999999
$ws.Range("N89").Value = -35545.75
$ws.Range("H105").Value = 500070.4
$ws.Range("I105").Value = 954801.9399999999
$ws.Range("K105").Value = 954801.9399999999
$ws.Range("M105").Value = -953054.9399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 753.1539
$ws.Range("I22").Value = 849.4
$ws.Range("K22").Value = 849.4
$ws.Range("M22").Value = -499.4
$ws.Range("H31").Value = 2596.439
$ws.Range("I31").Value = 2579.9355
$ws.Range("K31").Value = 2579.9355
$ws.Range("M31").Value = -2284.9355
$ws.Range("H34").Value = 2596.439
$ws.Range("I34").Value = 2579.9355
$ws.Range("K34").Value = 2579.9355
$ws.Range("M34").Value = -2377.9355
$ws.Range("H52").Value = 87249.25
$ws.Range("J52").Value = 99998.5
$ws.Range("L52").Value = 99998.5
$ws.Range("N52").Value = -100586.5
$ws.Range("H62").Value = 4245.077
$ws.Range("I62").Value = 4118
$ws.Range("J62").Value = 4944
$ws.Range("K62").Value = 4118
$ws.Range("L62").Value = 4944
$ws.Range("M62").Value = -3494
$ws.Range("N62").Value = -6192
$ws.Range("H65").Value = 4245.077
$ws.Range("I65").Value = 4118
$ws.Range("J65").Value = 4944
$ws.Range("K65").Value = 20590
$ws.Range("L65").Value = 24720
$ws.Range("M65").Value = -17470
$ws.Range("N65").Value = -30960
$ws.Range("H133").Value = 98315.2
$ws.Range("J133").Value = 98315.2
$ws.Range("L133").Value = 98315.2
$ws.Range("N133").Value = -103375.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 7341.4
$ws.Range("J60").Value = 16755
$ws.Range("L60").Value = 50265
$ws.Range("N60").Value = -50767

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2677759.5
$ws.Range("I132").Value = 4077.7144
$ws.Range("K132").Value = 12233.1432
$ws.Range("M132").Value = -9703.143199999999
$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("M137").Value = -105200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2626
$ws.Range("I61").Value = 2641.75
$ws.Range("K61").Value = 2641.75
$ws.Range("M61").Value = -2439.75
$ws.Range("H62").Value = 98666
$ws.Range("J62").Value = 98666
$ws.Range("L62").Value = 98666
$ws.Range("N62").Value = -99914
$ws.Range("H65").Value = 98666
$ws.Range("J65").Value = 98666
$ws.Range("L65").Value = 295998
$ws.Range("N65").Value = -302238
$ws.Range("H82").Value = 4740.174
$ws.Range("J82").Value = 8150.625
$ws.Range("L82").Value = 8150.625
$ws.Range("N82").Value = -8872.625
$ws.Range("H85").Value = 4740.174
$ws.Range("J85").Value = 8150.625
$ws.Range("L85").Value = 8150.625
$ws.Range("N85").Value = -10646.625
$ws.Range("H94").Value = 85164.5
$ws.Range("J94").Value = 85164.5
$ws.Range("L94").Value = 85164.5
$ws.Range("N94").Value = -86516.5
$ws.Range("H113").Value = 2626
$ws.Range("I113").Value = 2641.75
$ws.Range("K113").Value = 2641.75
$ws.Range("M113").Value = -471.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1259.3889
$ws.Range("I81").Value = 1039.3529
$ws.Range("K81").Value = 2078.7058
$ws.Range("M81").Value = -1017.7058
$ws.Range("H84").Value = 1259.3889
$ws.Range("I84").Value = 1039.3529
$ws.Range("K84").Value = 10393.529
$ws.Range("M84").Value = -5089.529
$ws.Range("H132").Value = 2577.2646
$ws.Range("I132").Value = 2501.45
$ws.Range("J132").Value = 2685.5715
$ws.Range("K132").Value = 7504.349999999999
$ws.Range("L132").Value = 8056.7145
$ws.Range("M132").Value = -4974.349999999999
$ws.Range("N132").Value = -13116.7145
$ws.Range("H136").Value = 7367.311
$ws.Range("I136").Value = 8364.4
$ws.Range("K136").Value = 25093.2
$ws.Range("M136").Value = -22543.2
